# Update the "last status check" timestamp in the header (F1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Last status check on: 25.02.2022 15:15"

# Tesco row (row 3): new price check -> shift Cena/Old Cena, recompute delta & old datum
$ws.Range("B3").Value = 36.9
$ws.Range("C3").Value = 36.7

# Delta Cena becomes a text label ("+0.2") instead of a plain number
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "+0.2"
$ws.Range("D3").Style = "Normal"

# Old Datum becomes a plain text timestamp instead of a formatted date serial
$ws.Range("E3").Value = "2022-02-25 15:17:15"
$ws.Range("E3").Style = "Normal"
